$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Einwohner/-innen pro m² Siedlungs- und Verkehrsfläche"
$ws.Range("B19").Value = "Millionen Einwohner/-innen"
$ws.Range("B22").Value = "Je 100 000 Einwohner/-innen"
$ws.Range("B23").Value = "Je 100 000 Einwohner/-innen unter 70 Jahren (ohne unter 1-Jährige)"
